$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Status")
$ws2 = $wb.Worksheets.Item("Waves&Experience")

# --- Status sheet: nerf enemy health (B column), bump Cruiser/CruiserElite
#     score (G column) and fix Laser Lv1 damage (D12) ---
$ws1.Range("D12").Value = 40    # Laser Lv1 damage 38 -> 40

$ws1.Range("B19").Value = 50    # Asteroid_medium health 60 -> 50
$ws1.Range("B20").Value = 150   # Asteroid_big health 180 -> 150
$ws1.Range("B21").Value = 50    # Razer health 60 -> 50
$ws1.Range("B22").Value = 100   # RazerBlood health 120 -> 100
$ws1.Range("B23").Value = 50    # Trident health 60 -> 50
$ws1.Range("B24").Value = 100   # TridentDark health 120 -> 100
$ws1.Range("B25").Value = 400   # Cruiser health 480 -> 400
$ws1.Range("G25").Value = 40    # Cruiser score 25 -> 40
$ws1.Range("B26").Value = 800   # CruiserElite health 960 -> 800
$ws1.Range("G26").Value = 60    # CruiserElite score 50 -> 60
$ws1.Range("B27").Value = 100   # Droid health 120 -> 100

# --- Waves&Experience sheet: wave 1-4 cost formula changed ---
$ws2.Range("C5").Formula = "=6*Status!G21 + 6*Status!G23"

# --- View / selection state ---
$ws2.Activate()
$ws2.Range("G19").Select()

$ws1.Activate()
$ws1.Range("H20").Select()
